$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("s__UBA11471 sp000434215-b-p")
$ws1.Rows("8:11").Delete()

$ws2 = $wb.Worksheets.Item("s__UBA11471 sp900542765-b-p")
$ws2.Rows("11:18").Delete()

$ws3 = $wb.Worksheets.Item("s__UBA11471 sp900547555-b-p")
$ws3.Rows("16:32").Delete()
